$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F5").Value = 0
$ws.Range("F8").Value = 2
$ws.Range("F9").Value = 0
$ws.Range("F15").Value = 0
$ws.Range("F17").Value = -12
$ws.Range("F19").Value = -8
